$wb = $excel.ActiveWorkbook

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2100
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 2100
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 2100
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -2668

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14162.6
$ws.Range("I32").Value = 12402.889
$ws.Range("J32").Value = 30000
$ws.Range("K32").Value = 12402.889
$ws.Range("L32").Value = 30000
$ws.Range("M32").Value = -12115.889
$ws.Range("N32").Value = -30574

# ARM row 41
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 26078
$ws.Range("I41").Value = 2156
$ws.Range("K41").Value = 2156
$ws.Range("M41").Value = -1742

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4962.6665
$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 5000
$ws.Range("N102").Value = -8244

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2509.375
$ws.Range("J122").Value = 2796.3333
$ws.Range("L122").Value = 8388.999899999999
$ws.Range("N122").Value = -13288.9999

# BSM row 6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 62333
$ws.Range("I6").Value = 62333
$ws.Range("K6").Value = 62333
$ws.Range("M6").Value = -62220

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2640.8572
$ws.Range("I20").Value = 2005.6666
$ws.Range("K20").Value = 2005.6666
$ws.Range("M20").Value = -1758.6666

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1688.3704
$ws.Range("I86").Value = 1051.762
$ws.Range("J86").Value = 3916.5
$ws.Range("K86").Value = 1051.762
$ws.Range("L86").Value = 3916.5
$ws.Range("M86").Value = 71.23800000000006
$ws.Range("N86").Value = -6162.5

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1688.3704
$ws.Range("I89").Value = 1051.762
$ws.Range("J89").Value = 3916.5
$ws.Range("K89").Value = 5258.809999999999
$ws.Range("L89").Value = 19582.5
$ws.Range("M89").Value = 357.1900000000005
$ws.Range("N89").Value = -30814.5

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4775
$ws.Range("I99").Value = 4775
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4775
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3277
$ws.Range("N99").ClearContents()

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1867.4615
$ws.Range("J105").Value = 1884
$ws.Range("L105").Value = 1884
$ws.Range("N105").Value = -5378

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1296.6666
$ws.Range("I107").Value = 1296.6666
$ws.Range("K107").Value = 1296.6666
$ws.Range("M107").Value = 623.3334

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2590.7896
$ws.Range("I16").Value = 1240.1428
$ws.Range("K16").Value = 1240.1428
$ws.Range("M16").Value = -953.1428000000001

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 251.625
$ws.Range("I22").Value = 222.2
$ws.Range("J22").Value = 300.66666
$ws.Range("K22").Value = 222.2
$ws.Range("L22").Value = 300.66666
$ws.Range("M22").Value = 127.8
$ws.Range("N22").Value = -1000.66666

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6543.5
$ws.Range("I99").Value = 6692.5713
$ws.Range("J99").Value = 5500
$ws.Range("K99").Value = 6692.5713
$ws.Range("L99").Value = 5500
$ws.Range("M99").Value = -5194.5713
$ws.Range("N99").Value = -8496

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 620.25
$ws.Range("I107").Value = 750
$ws.Range("K107").Value = 750
$ws.Range("M107").Value = 1170

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2590.7896
$ws.Range("I113").Value = 1240.1428
$ws.Range("K113").Value = 1240.1428
$ws.Range("M113").Value = 929.8571999999999

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6543.5
$ws.Range("I126").Value = 6692.5713
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 20077.7139
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -17607.7139
$ws.Range("N126").Value = -21440

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1309.6666
$ws.Range("I132").Value = 1183.8572
$ws.Range("J132").Value = 1750
$ws.Range("K132").Value = 3551.5716
$ws.Range("L132").Value = 5250
$ws.Range("M132").Value = -1021.5716
$ws.Range("N132").Value = -10310

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 253.75
$ws.Range("I12").Value = 240
$ws.Range("J12").Value = 267.5
$ws.Range("K12").Value = 720
$ws.Range("L12").Value = 802.5
$ws.Range("M12").Value = -547
$ws.Range("N12").Value = -1148.5

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 125
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2919.5312
$ws.Range("J80").Value = 3017.6072
$ws.Range("L80").Value = 3017.6072
$ws.Range("N80").Value = -5013.6072

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2919.5312
$ws.Range("J83").Value = 3017.6072
$ws.Range("L83").Value = 15088.036
$ws.Range("N83").Value = -25072.036

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1424.875
$ws.Range("I102").Value = 879.8
$ws.Range("K102").Value = 879.8
$ws.Range("M102").Value = 742.2

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1900
$ws.Range("I113").Value = 1800
$ws.Range("K113").Value = 1800
$ws.Range("M113").Value = 370

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1865.6666
$ws.Range("I132").Value = 1838.8
$ws.Range("K132").Value = 5516.4
$ws.Range("M132").Value = -2986.4

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 13074.25
$ws.Range("I16").Value = 12999
$ws.Range("J16").Value = 13149.5
$ws.Range("K16").Value = 12999
$ws.Range("L16").Value = 13149.5
$ws.Range("M16").Value = -12829
$ws.Range("N16").Value = -13489.5

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1985.5714
$ws.Range("J22").Value = 1999.8
$ws.Range("L22").Value = 1999.8
$ws.Range("N22").Value = -2589.8

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1985.5714
$ws.Range("J27").Value = 1999.8
$ws.Range("L27").Value = 1999.8
$ws.Range("N27").Value = -2213.8

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8812
$ws.Range("I40").Value = 8614
$ws.Range("K40").Value = 8614
$ws.Range("M40").Value = -8478

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4161.263
$ws.Range("I61").Value = 3474.353
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 3474.353
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -3272.353
$ws.Range("N61").Value = -10404

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 361
$ws.Range("I93").Value = 361
$ws.Range("K93").Value = 361
$ws.Range("M93").Value = 887

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4161.263
$ws.Range("I113").Value = 3474.353
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 3474.353
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -1304.353
$ws.Range("N113").Value = -14340

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1080
$ws.Range("N107").ClearContents()

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1855.3684
$ws.Range("I136").Value = 1826.5883
$ws.Range("K136").Value = 5479.7649
$ws.Range("M136").Value = -2929.7649
